$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-10, columns B-G (B=ME, C=MAE, D=MSE, E=RMSE, F=SE, G=N)
$data = @{
    2  = @(0.1417620482754679, 0.5996277706543274, 0.8169135268031753, 0.9038326873947276, 0.9017087228643812, 50)
    3  = @(0.09584318641535604, 0.657502611948801, 1.03928723612716, 1.019454381582207, 1.025456856189602, 49)
    4  = @(0.126342760679196, 0.5716138419826224, 0.7986488544558821, 0.8936715584910835, 0.894057751257528, 48)
    5  = @(0.132571133247816, 0.6727079862838818, 1.180166600671597, 1.086354730588309, 1.089892299018441, 47)
    6  = @(0.1115124553593621, 0.6125314060269084, 0.941813663750834, 0.9704708464198365, 0.9746955907510105, 46)
    7  = @(0.0971729415500782, 0.5646146643790109, 0.6953729178386839, 0.8338902312886775, 0.8406640857472294, 34)
    8  = @(0.1456103435239261, 0.5758738547837706, 0.7130976704205472, 0.8444511059975865, 0.8446993703979081, 33)
    9  = @(-0.01600337040817491, 0.5629050319730056, 0.5311176513683245, 0.7287781907880644, 0.7524973840956725, 16)
    10 = @(-0.08801024940629837, 0.472497539938865, 0.4690803580242086, 0.6848944137779257, 0.7159566669255664, 10)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
